$wb = $excel.ActiveWorkbook

# Cell value updates scraped from the latest market-board price refresh.
# Each sheet entry is a list of (cellRef, newValue) pairs; a $null value
# clears the cell entirely (matches source rows where Profit columns are
# omitted rather than written as 0).

$updates_ALC = @(
    @("H28", 942.625),
    @("J28", 1400.2222),
    @("L28", 1400.2222),
    @("N28", -2370.2222),
    @("H61", 558.8570999999999),
    @("I61", 649.1667),
    @("J61", 17),
    @("K61", 1947.5001),
    @("L61", 51),
    @("M61", -1775.5001),
    @("N61", -395),
    @("H88", 1073.5),
    @("I88", 0),
    @("J88", 1073.5),
    @("K88", 0),
    @("L88", 1073.5),
    @("M88", $null),
    @("N88", -1885.5),
    @("H91", 1073.5),
    @("I91", 0),
    @("J91", 1073.5),
    @("K91", 0),
    @("L91", 1073.5),
    @("M91", $null),
    @("N91", -3881.5),
    @("H96", 11364217),
    @("I96", 27778314),
    @("J96", 612.1539),
    @("K96", 83334942),
    @("L96", 1836.4617),
    @("M96", -83333569),
    @("N96", -4582.4617),
    @("H103", 147),
    @("J103", 163.71428),
    @("L103", 491.14284),
    @("N103", -1663.14284),
    @("H106", 7578283.5),
    @("I106", 13890632),
    @("K106", 13890632),
    @("M106", -13890001),
    @("H116", 6728.1816),
    @("I116", 3502.5),
    @("J116", 7445),
    @("K116", 3502.5),
    @("L116", 7445),
    @("M116", -60.5),
    @("N116", -14329),
    @("H129", 213868.42),
    @("J129", 228418.55),
    @("L129", 685255.6499999999),
    @("N129", -695255.6499999999),
    @("H141", 5133.3335),
    @("I141", 0),
    @("J141", 5133.3335),
    @("K141", 0),
    @("L141", 15400.0005),
    @("M141", $null),
    @("N141", -25760.0005),
)

$updates_ARM = @(
    @("H5", 0),
    @("I5", 0),
    @("J5", 0),
    @("K5", 0),
    @("L5", 0),
    @("M5", $null),
    @("N5", $null),
    @("H32", 6534.7363),
    @("I32", 4401.4263),
    @("J32", 18364.908),
    @("K32", 4401.4263),
    @("L32", 18364.908),
    @("M32", -4114.4263),
    @("N32", -18938.908),
    @("H74", 38463020),
    @("I74", 71429100),
    @("K74", 71429100),
    @("M74", -71428226),
    @("H77", 38463020),
    @("I77", 71429100),
    @("K77", 357145500),
    @("M77", -357141132),
    @("H88", 126511.5),
    @("I88", 1408.2),
    @("J88", 335017),
    @("K88", 1408.2),
    @("L88", 335017),
    @("M88", -1002.2),
    @("N88", -335829),
    @("H91", 126511.5),
    @("I91", 1408.2),
    @("J91", 335017),
    @("K91", 1408.2),
    @("L91", 335017),
    @("M91", -4.200000000000045),
    @("N91", -337825),
)

$updates_BSM = @(
    @("H4", 0),
    @("I4", 0),
    @("J4", 0),
    @("K4", 0),
    @("L4", 0),
    @("M4", $null),
    @("N4", $null),
    @("H105", 2779362),
    @("I105", 1419),
    @("J105", 8335248.5),
    @("K105", 1419),
    @("L105", 8335248.5),
    @("M105", 328),
    @("N105", -8338742.5),
)

$updates_CRP = @(
    @("H31", 3463.7715),
    @("I31", 2462.389),
    @("J31", 4524.0586),
    @("K31", 2462.389),
    @("L31", 4524.0586),
    @("M31", -2167.389),
    @("N31", -5114.0586),
    @("H34", 3463.7715),
    @("I34", 2462.389),
    @("J34", 4524.0586),
    @("K34", 2462.389),
    @("L34", 4524.0586),
    @("M34", -2260.389),
    @("N34", -4928.0586),
    @("H86", 19986.5),
    @("I86", 3680),
    @("J86", 47164),
    @("K86", 3680),
    @("L86", 47164),
    @("M86", -2557),
    @("N86", -49410),
    @("H89", 19986.5),
    @("I89", 3680),
    @("J89", 47164),
    @("K89", 18400),
    @("L89", 235820),
    @("M89", -12784),
    @("N89", -247052),
    @("H107", 1291.1305),
    @("I107", 583.1111),
    @("J107", 1746.2858),
    @("K107", 583.1111),
    @("L107", 1746.2858),
    @("M107", 1336.8889),
    @("N107", -5586.2858),
)

$updates_CUL = @(
    @("H97", 401),
    @("I97", 266),
    @("J97", 536),
    @("K97", 798),
    @("L97", 1608),
    @("M97", -302),
    @("N97", -2600),
    @("H131", 719.8099999999999),
    @("I131", 476.66666),
    @("J131", 727.3299),
    @("K131", 1429.99998),
    @("L131", 2181.9897),
    @("M131", 3610.00002),
    @("N131", -12261.9897),
)

$updates_GSM = @(
    @("H80", 4203.2),
    @("I80", 3299.75),
    @("J80", 5235.7144),
    @("K80", 3299.75),
    @("L80", 5235.7144),
    @("M80", -2301.75),
    @("N80", -7231.7144),
    @("H83", 4203.2),
    @("I83", 3299.75),
    @("J83", 5235.7144),
    @("K83", 16498.75),
    @("L83", 26178.572),
    @("M83", -11506.75),
    @("N83", -36162.572),
    @("H95", 21007.75),
    @("J95", 21007.75),
    @("L95", 21007.75),
    @("N95", -26499.75),
    @("H97", 2152.647),
    @("I97", 2144.5),
    @("J97", 2164.2856),
    @("K97", 2144.5),
    @("L97", 2164.2856),
    @("M97", -1648.5),
    @("N97", -3156.2856),
    @("H132", 34149.312),
    @("I132", 4236.5),
    @("J132", 44120.25),
    @("K132", 12709.5),
    @("L132", 132360.75),
    @("M132", -10179.5),
    @("N132", -137420.75),
)

$updates_LTW = @(
    @("H7", 6063.2104),
    @("I7", 4125.5),
    @("J7", 7472.4546),
    @("K7", 4125.5),
    @("L7", 7472.4546),
    @("M7", -4013.5),
    @("N7", -7696.4546),
    @("H126", 6063.2104),
    @("I126", 4125.5),
    @("J126", 7472.4546),
    @("K126", 12376.5),
    @("L126", 22417.3638),
    @("M126", -9906.5),
    @("N126", -27357.3638),
)

$updates_WVR = @(
    @("H126", 1600.2222),
    @("I126", 1650.25),
    @("J126", 1200),
    @("K126", 4950.75),
    @("L126", 3600),
    @("M126", -2480.75),
    @("N126", -8540),
)

$sheetUpdates = @{
    "ALC" = $updates_ALC
    "ARM" = $updates_ARM
    "BSM" = $updates_BSM
    "CRP" = $updates_CRP
    "CUL" = $updates_CUL
    "GSM" = $updates_GSM
    "LTW" = $updates_LTW
    "WVR" = $updates_WVR
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $updates = $sheetUpdates[$sheetName]
    foreach ($u in $updates) {
        $cellRef = $u[0]
        $newVal = $u[1]
        if ($null -eq $newVal) {
            $ws.Range($cellRef).Value = ""
        } else {
            $ws.Range($cellRef).Value = $newVal
        }
    }
}
